$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure the target cells are formatted as Text so Excel stores the
# values as plain strings (inline strings) instead of converting them
# into numbers or dates.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("I2").NumberFormat = "@"

$ws.Range("C2").Value = "ALEJANDRA S FASHIÓON INC"
$ws.Range("F2").Value = "35,783.00"
$ws.Range("H2").Value = "5293"
$ws.Range("I2").Value = "2024-11-08"
